$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 14 (relabeled scan + refreshed values) ---
$ws.Range("B14").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C14").Value = 0.999519819632665
$ws.Range("D14").Value = 1.017562233915762
$ws.Range("E14").Value = 0.9828356523002032
$ws.Range("F14").Value = 0.9922549711631075
$ws.Range("G14").Value = 0.999519819632665
$ws.Range("H14").Value = 1.017562233915762
$ws.Range("I14").Value = 0.9834194394324258
$ws.Range("J14").Value = 0.9833850312650245
$ws.Range("K14").Value = 0.988071276157061
$ws.Range("L14").Value = 1.002057048397828
$ws.Range("M14").Value = 0.9994806574392824
$ws.Range("N14").Value = 1.000198943107983
$ws.Range("O14").Value = 1.004908602539435
$ws.Range("P14").Value = 0.9999725686162101
$ws.Range("Q14").Value = 0.9975509524596907
$ws.Range("R14").Value = 0.9999725686162101
$ws.Range("S14").Value = 0.9980431692529343
$ws.Range("T14").Value = 0.9983384993288805
$ws.Range("U14").Value = 0.9936381840330096

# --- Append two new rows (22, 23 index / sheet rows 24-25), matching
#     the existing table's formatting by copying row 23's formats down ---
$ws.Range("A23:U23").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 24
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "RotRing Axis-Y Res-5.0 Theta-2.5 "
$ws.Range("C24").Value = 1.064566973029285
$ws.Range("D24").Value = 1.026181214278875
$ws.Range("E24").Value = 0.9493745435035073
$ws.Range("F24").Value = 0.9960352639173311
$ws.Range("G24").Value = 1.064566973029285
$ws.Range("H24").Value = 1.026181214278875
$ws.Range("I24").Value = 0.9728871906349001
$ws.Range("J24").Value = 0.9597549947074582
$ws.Range("K24").Value = 1.018186068794158
$ws.Range("L24").Value = 0.993989054773298
$ws.Range("M24").Value = 1.064560347927988
$ws.Range("N24").Value = 0.9877778788911912
$ws.Range("O24").Value = 1.011108239098103
$ws.Range("P24").Value = 1.013374243603889
$ws.Range("Q24").Value = 0.9905303405665711
$ws.Range("R24").Value = 1.013374243603889
$ws.Range("S24").Value = 1.00903949868225
$ws.Range("T24").Value = 1.020144993551657
$ws.Range("U24").Value = 0.9976219129548516

# Row 25
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "RotRing Axis-Y Res-5.0 Theta-5.0 "
$ws.Range("C25").Value = 1.055669658786997
$ws.Range("D25").Value = 1.016021296484453
$ws.Range("E25").Value = 0.9522480769540386
$ws.Range("F25").Value = 0.9937391728448696
$ws.Range("G25").Value = 1.055669658786997
$ws.Range("H25").Value = 1.016021296484453
$ws.Range("I25").Value = 0.9768538994046883
$ws.Range("J25").Value = 0.9639748353499836
$ws.Range("K25").Value = 1.012134108143645
$ws.Range("L25").Value = 0.9933785524706517
$ws.Range("M25").Value = 1.055658804731599
$ws.Range("N25").Value = 0.9841346867192458
$ws.Range("O25").Value = 1.004880234664661
$ws.Range("P25").Value = 1.007979677408496
$ws.Range("Q25").Value = 0.9873361820944537
$ws.Range("R25").Value = 1.007979677408496
$ws.Range("S25").Value = 1.00441955126759
$ws.Range("T25").Value = 1.014669572771471
$ws.Range("U25").Value = 0.995502450054916
